$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("revenue")
$ws.Range("B5").Value = 1
